$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the leadership bio links to point to individual pages instead of
# the generic Leadership_team.html page.
$ws.Range("F2").Value = "[Lauren Chenarides](https://dataifa.github.io/difa-project/lauren_chenarides.html), [Drew Hanks](https://dataifa.github.io/difa-project/drew_hanks.html)"
$ws.Range("F10").Value = "[Lauren Chenarides](https://dataifa.github.io/difa-project/lauren_chenarides.html), [Drew Hanks](https://dataifa.github.io/difa-project/drew_hanks.html)"
$ws.Range("F8").Value = "[Andi Carlson](https://dataifa.github.io/difa-project/andi_carlson.html)"

# Update the active selection to reflect where the author left off editing.
$ws.Range("F10").Select()

$wb.Save()
